$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.706.34"
$ws.Range("E2").Value = "  +1.76%  "

$ws.Range("D3").Value = "2.523.73"
$ws.Range("E3").Value = "  -1.62%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.62%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.529"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.41%  "

$ws.Range("D9").Value = "2.523.60"
$ws.Range("E9").Value = "  -1.62%  "

$ws.Range("E10").Value = "  +1.04%  "

$ws.Range("E11").Value = "  +1.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.343"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.21%  "

$ws.Range("D15").Value = "2.987.31"
$ws.Range("E15").Value = "  -1.78%  "

$ws.Range("E16").Value = "  -0.05%  "

$ws.Range("D17").Value = "67.600.41"
$ws.Range("E17").Value = "  +1.81%  "

$ws.Range("D18").Value = "2.536.75"
$ws.Range("E18").Value = "  -0.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.62%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "357.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.61%  "

$ws.Range("E25").Value = "  +0.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "70.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.95%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.657.74"
$ws.Range("E28").Value = "  -1.82%  "

$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.67%  "

$ws.Range("D30").Value = "0.0₃0984"
$ws.Range("E30").Value = "  -0.54%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "553.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.34%  "

$ws.Range("E34").Value = "  +1.49%  "

$ws.Range("E35").Value = "  -0.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.72"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.43%  "

$ws.Range("E41").Value = "  +2.22%  "

$ws.Range("E42").Value = "  -0.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.94%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.97%  "

$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "147.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.85%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.558"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.10%  "

$ws.Range("D48").Value = "0.0₆0277"
$ws.Range("E48").Value = "  -2.88%  "

$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0757"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.31%  "

